$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Status column (M) changes from "A" to "B" for rows 2, 3 and 4 (row 5 keeps "A")
$ws.Range("M2").Value = "B"
$ws.Range("M3").Value = "B"
$ws.Range("M4").Value = "B"

# Saldo (N) for row 5 is fully settled -> 0
$ws.Range("N5").Value = 0.0

# Column N (Saldo) narrows now that its widest value is "0.0" instead of "1565.01"
# (closest attainable width given the engine's character-width rounding)
$ws.Columns.Item(14).ColumnWidth = 5.14
